$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stock symbol cells that changed (columns B/C/D/E/F, rows 2-42)
$ws.Range("B2").Value = "NSE:AARTIDRUGS"
$ws.Range("C2").Value = "NSE:ABBOTINDIA"
$ws.Range("F2").Value = "NSE:HINDCOPPER"
$ws.Range("B3").Value = "NSE:ALANKIT"
$ws.Range("C3").Value = "NSE:BSE"
$ws.Range("F3").Value = "NSE:IRFC"
$ws.Range("B4").Value = "NSE:ALMONDZ"
$ws.Range("C4").Value = "NSE:FAZE3Q"
$ws.Range("B5").Value = "NSE:ANIKINDS"
$ws.Range("C5").Value = "NSE:JINDWORLD"
$ws.Range("B6").Value = "NSE:ARVEE"
$ws.Range("B7").Value = "NSE:BEML"
$ws.Range("B8").Value = "NSE:CAPTRUST"
$ws.Range("B9").Value = "NSE:CENTENKA"
$ws.Range("B10").Value = "NSE:CYIENTDLM"
$ws.Range("B11").Value = "NSE:DCI"
$ws.Range("B12").Value = "NSE:DICIND"
$ws.Range("B13").Value = "NSE:ENDURANCE"
$ws.Range("B14").Value = "NSE:EPIGRAL"
$ws.Range("B15").Value = "NSE:GENUSPAPER"
$ws.Range("B16").Value = "NSE:GLOBUSSPR"
$ws.Range("B17").Value = "NSE:GOYALALUM"
$ws.Range("B18").Value = "NSE:GREENPOWER"
$ws.Range("B19").Value = "NSE:HATHWAY"
$ws.Range("B20").Value = "NSE:HIMATSEIDE"
$ws.Range("B21").Value = "NSE:HINDCOPPER"
$ws.Range("B22").Value = "NSE:HMAAGRO"
$ws.Range("B23").Value = "NSE:INDOAMIN"
$ws.Range("B24").Value = "NSE:JAIBALAJI"
$ws.Range("B25").Value = "NSE:JASH"
$ws.Range("B26").Value = "NSE:JAYBARMARU"
$ws.Range("B27").Value = "NSE:KCPSUGIND"
$ws.Range("B28").Value = "NSE:KIOCL"
$ws.Range("B29").Value = "NSE:KIRIINDUS"
$ws.Range("B30").Value = "NSE:KPIL"
$ws.Range("B31").Value = "NSE:LMW"
$ws.Range("B32").Value = "NSE:MAHLIFE"
$ws.Range("B33").Value = "NSE:MIDHANI"
$ws.Range("B34").Value = "NSE:NETWORK18"
$ws.Range("B35").Value = "NSE:POWERMECH"
$ws.Range("B36").Value = "NSE:PRITIKAUTO"
$ws.Range("B37").Value = "NSE:RAMCOIND"
$ws.Range("B38").Value = "NSE:RANASUG"
$ws.Range("B39").Value = "NSE:RELIANCE"
$ws.Range("B40").Value = "NSE:RSWM"
$ws.Range("B41").Value = "NSE:RTNINDIA"
$ws.Range("B42").Value = "NSE:SADBHIN"

# Clear cells that became empty
$ws.Range("D2").ClearContents()
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("D5").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("F6").ClearContents()
$ws.Range("F7").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("F9").ClearContents()

# Remove the now-unused trailing rows (43 and 44)
$ws.Rows("43:44").Delete()
